$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the per-category detail rows to include the category name as a
# prefix (e.g. "     New nominations" -> "     Civilian, New nominations").

# Civilian
$ws.Range("A7").Value  = "     Civilian, New nominations"
$ws.Range("A8").Value  = "     Civilian, Confirmed "
$ws.Range("A9").Value  = "     Civilian, Unconfirmed "
$ws.Range("A10").Value = "     Civilian, Withdrawn "
$ws.Range("A11").Value = "     Civilian, Returned to White House "

# Other Civilian
$ws.Range("A13").Value = "     Other Civilian, New nominations"
$ws.Range("A14").Value = "     Other Civilian, Confirmed "
$ws.Range("A15").Value = "     Other Civilian, Unconfirmed "
$ws.Range("A16").Value = "     Other Civilian, Returned to White House "

# Air Force
$ws.Range("A18").Value = "     Air Force, New nominations"
$ws.Range("A19").Value = "     Air Force, Confirmed "
$ws.Range("A20").Value = "     Air Force, Unconfirmed "

# Army
$ws.Range("A22").Value = "     Army, New nominations"
$ws.Range("A23").Value = "     Army, Confirmed "
$ws.Range("A24").Value = "     Army, Unconfirmed "

# Navy
$ws.Range("A26").Value = "     Navy, New nominations"
$ws.Range("A27").Value = "     Navy, Confirmed "
$ws.Range("A28").Value = "     Navy, Unconfirmed "

# Marine Corps
$ws.Range("A30").Value = "     Marine Corps, New nominations"
$ws.Range("A31").Value = "     Marine Corps, Confirmed "
$ws.Range("A32").Value = "     Marine Corps, Unconfirmed "

# The "Summary" section header (row 33) is removed and its former value
# row 34 ("Total nominations received this session ") is renamed/merged
# into row 33 as "Total new nominations", keeping its value (23640). That
# row previously had no value cell, so pick up the same "#,##0" number
# format the other big-total cells (e.g. B34) use by copying formats only.
$ws.Range("A33").Value = "Total new nominations"
$ws.Range("B33").Value = 23640
$ws.Range("B34").Copy()
$ws.Range("B33").PasteSpecial(-4122)

$ws.Range("A34").Value = "Total confirmed "
$ws.Range("B34").Value = 22468

$ws.Range("A35").Value = "Total unconfirmed "
$ws.Range("B35").Value = 1150

# B36 previously held the "#,##0"-formatted 1150 total; the new content is
# a small plain-number count, so pick up B37's (General) number format.
$ws.Range("A36").Value = "Total withdrawn "
$ws.Range("B36").Value = 12
$ws.Range("B37").Copy()
$ws.Range("B36").PasteSpecial(-4122)

$ws.Range("A37").Value = "Total returned"
$ws.Range("B37").Value = 10

# The old row 38 ("Total Returned to White House ") is now redundant
# (its value moved up into row 37) so delete the trailing row.
$ws.Rows("38").Delete()
